$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (price_level) - shifts old D (rating) to E, old E (user_ratings_total) to F
$ws.Columns("D:D").Insert()
$ws.Cells.Item(1,4).Value = "price_level"

# Clear the existing data rows (2-59); the table is being refreshed from a new data pull
# that adds/removes a couple of venues and re-sorts everything by name.
$ws.Range("A2:F59").ClearContents()

# Re-populate every row of this fresh, alphabetically-sorted dataset.
$ws.Cells.Item(2,1).Value = 19
$ws.Cells.Item(2,2).Value = "OPERATIONAL"
$ws.Cells.Item(2,3).Value = "A Classic Expo Design"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1

$ws.Cells.Item(3,1).Value = 28
$ws.Cells.Item(3,2).Value = "OPERATIONAL"
$ws.Cells.Item(3,3).Value = "Aesthetically Sew"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0

$ws.Cells.Item(4,1).Value = 17
$ws.Cells.Item(4,2).Value = "OPERATIONAL"
$ws.Cells.Item(4,3).Value = "AgTech Expo"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0

$ws.Cells.Item(5,1).Value = 38
$ws.Cells.Item(5,2).Value = "OPERATIONAL"
$ws.Cells.Item(5,3).Value = "American Rodent Supply"
$ws.Cells.Item(5,5).Value = 4.6
$ws.Cells.Item(5,6).Value = 28

$ws.Cells.Item(6,1).Value = 44
$ws.Cells.Item(6,2).Value = "OPERATIONAL"
$ws.Cells.Item(6,3).Value = "Biltwell Event Center"
$ws.Cells.Item(6,5).Value = 4.7
$ws.Cells.Item(6,6).Value = 373

$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "OPERATIONAL"
$ws.Cells.Item(7,3).Value = "Cabinet and Stone Expo"
$ws.Cells.Item(7,5).Value = 4.6
$ws.Cells.Item(7,6).Value = 31

$ws.Cells.Item(8,1).Value = 43
$ws.Cells.Item(8,2).Value = "OPERATIONAL"
$ws.Cells.Item(8,3).Value = "Circle City Coatings"
$ws.Cells.Item(8,5).Value = 4.6
$ws.Cells.Item(8,6).Value = 19

$ws.Cells.Item(9,1).Value = 29
$ws.Cells.Item(9,2).Value = "OPERATIONAL"
$ws.Cells.Item(9,3).Value = "Convention center"
$ws.Cells.Item(9,5).Value = 4.6
$ws.Cells.Item(9,6).Value = 473

$ws.Cells.Item(10,1).Value = 25
$ws.Cells.Item(10,2).Value = "OPERATIONAL"
$ws.Cells.Item(10,3).Value = "Core & Main"
$ws.Cells.Item(10,5).Value = 4.3
$ws.Cells.Item(10,6).Value = 3

$ws.Cells.Item(11,1).Value = 0
$ws.Cells.Item(11,2).Value = "OPERATIONAL"
$ws.Cells.Item(11,3).Value = "Creation Evidence Expo"
$ws.Cells.Item(11,5).Value = 4.8
$ws.Cells.Item(11,6).Value = 5

$ws.Cells.Item(12,1).Value = 51
$ws.Cells.Item(12,2).Value = "OPERATIONAL"
$ws.Cells.Item(12,3).Value = "Curtain Call Dance Center"
$ws.Cells.Item(12,5).Value = 4.4
$ws.Cells.Item(12,6).Value = 9

$ws.Cells.Item(13,1).Value = 4
$ws.Cells.Item(13,2).Value = "OPERATIONAL"
$ws.Cells.Item(13,3).Value = "Curvature Expo"
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 0

$ws.Cells.Item(14,1).Value = 15
$ws.Cells.Item(14,2).Value = "OPERATIONAL"
$ws.Cells.Item(14,3).Value = "Don Mitchell Pro Shop"
$ws.Cells.Item(14,5).Value = 3.3
$ws.Cells.Item(14,6).Value = 7

$ws.Cells.Item(15,1).Value = 42
$ws.Cells.Item(15,2).Value = "OPERATIONAL"
$ws.Cells.Item(15,3).Value = "Don Mitchell Pro Shops"
$ws.Cells.Item(15,5).Value = 4.4
$ws.Cells.Item(15,6).Value = 16

$ws.Cells.Item(16,1).Value = 49
$ws.Cells.Item(16,2).Value = "OPERATIONAL"
$ws.Cells.Item(16,3).Value = "Elegant Stylez"
$ws.Cells.Item(16,5).Value = 3.1
$ws.Cells.Item(16,6).Value = 3696

$ws.Cells.Item(17,1).Value = 57
$ws.Cells.Item(17,2).Value = "OPERATIONAL"
$ws.Cells.Item(17,3).Value = "Elements Financial Pavilion"
$ws.Cells.Item(17,5).Value = 4.3
$ws.Cells.Item(17,6).Value = 87

$ws.Cells.Item(18,1).Value = 2
$ws.Cells.Item(18,2).Value = "OPERATIONAL"
$ws.Cells.Item(18,3).Value = "Expo Arts"
$ws.Cells.Item(18,5).Value = 2.5
$ws.Cells.Item(18,6).Value = 4

$ws.Cells.Item(19,1).Value = 41
$ws.Cells.Item(19,2).Value = "OPERATIONAL"
$ws.Cells.Item(19,3).Value = "FYE"
$ws.Cells.Item(19,4).Value = 2
$ws.Cells.Item(19,5).Value = 4.1
$ws.Cells.Item(19,6).Value = 638

$ws.Cells.Item(20,1).Value = 27
$ws.Cells.Item(20,2).Value = "OPERATIONAL"
$ws.Cells.Item(20,3).Value = "Fern"
$ws.Cells.Item(20,5).Value = 4.3
$ws.Cells.Item(20,6).Value = 34

$ws.Cells.Item(21,1).Value = 47
$ws.Cells.Item(21,2).Value = "OPERATIONAL"
$ws.Cells.Item(21,3).Value = "Floor & Decor"
$ws.Cells.Item(21,5).Value = 4.4
$ws.Cells.Item(21,6).Value = 281

$ws.Cells.Item(22,1).Value = 34
$ws.Cells.Item(22,2).Value = "OPERATIONAL"
$ws.Cells.Item(22,3).Value = "Food Specialties Inc"
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0

$ws.Cells.Item(23,1).Value = 52
$ws.Cells.Item(23,2).Value = "OPERATIONAL"
$ws.Cells.Item(23,3).Value = "Great Day Tattoo"
$ws.Cells.Item(23,5).Value = 5
$ws.Cells.Item(23,6).Value = 18

$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "OPERATIONAL"
$ws.Cells.Item(24,3).Value = "Harvest Pavillion"
$ws.Cells.Item(24,5).Value = 4.6
$ws.Cells.Item(24,6).Value = 18

$ws.Cells.Item(25,1).Value = 33
$ws.Cells.Item(25,2).Value = "OPERATIONAL"
$ws.Cells.Item(25,3).Value = "Hoosier Lottery Hall"
$ws.Cells.Item(25,5).Value = 5
$ws.Cells.Item(25,6).Value = 3

$ws.Cells.Item(26,1).Value = 55
$ws.Cells.Item(26,2).Value = "OPERATIONAL"
$ws.Cells.Item(26,3).Value = "Hoosier Trim Products"
$ws.Cells.Item(26,5).Value = 4.7
$ws.Cells.Item(26,6).Value = 3

$ws.Cells.Item(27,1).Value = 53
$ws.Cells.Item(27,2).Value = "OPERATIONAL"
$ws.Cells.Item(27,3).Value = "INKSTINCT TATTOO"
$ws.Cells.Item(27,5).Value = 4.8
$ws.Cells.Item(27,6).Value = 116

$ws.Cells.Item(28,1).Value = 3
$ws.Cells.Item(28,2).Value = "OPERATIONAL"
$ws.Cells.Item(28,3).Value = "Indiana Black Expo Inc"
$ws.Cells.Item(28,5).Value = 4.3
$ws.Cells.Item(28,6).Value = 39

$ws.Cells.Item(29,1).Value = 24
$ws.Cells.Item(29,2).Value = "OPERATIONAL"
$ws.Cells.Item(29,3).Value = "Indiana Convention Center"
$ws.Cells.Item(29,5).Value = 4.5
$ws.Cells.Item(29,6).Value = 528

$ws.Cells.Item(30,1).Value = 11
$ws.Cells.Item(30,2).Value = "OPERATIONAL"
$ws.Cells.Item(30,3).Value = "Indiana Fishing Expo"
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 0

$ws.Cells.Item(31,1).Value = 23
$ws.Cells.Item(31,2).Value = "OPERATIONAL"
$ws.Cells.Item(31,3).Value = "Indiana Flower & Patio Show"
$ws.Cells.Item(31,5).Value = 4.3
$ws.Cells.Item(31,6).Value = 60

$ws.Cells.Item(32,1).Value = 5
$ws.Cells.Item(32,2).Value = "OPERATIONAL"
$ws.Cells.Item(32,3).Value = "Indiana Latino Expo"
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(32,6).Value = 0

$ws.Cells.Item(33,1).Value = 20
$ws.Cells.Item(33,2).Value = "OPERATIONAL"
$ws.Cells.Item(33,3).Value = "Indiana State Fairgrounds & Event Center"
$ws.Cells.Item(33,5).Value = 4.4
$ws.Cells.Item(33,6).Value = 1344

$ws.Cells.Item(34,1).Value = 40
$ws.Cells.Item(34,2).Value = "OPERATIONAL"
$ws.Cells.Item(34,3).Value = "Indiana State Numismatic Association"
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(34,6).Value = 0

$ws.Cells.Item(35,1).Value = 30
$ws.Cells.Item(35,2).Value = "OPERATIONAL"
$ws.Cells.Item(35,3).Value = "Indianapolis Auto Show"
$ws.Cells.Item(35,5).Value = 3.1
$ws.Cells.Item(35,6).Value = 51

$ws.Cells.Item(36,1).Value = 14
$ws.Cells.Item(36,2).Value = "OPERATIONAL"
$ws.Cells.Item(36,3).Value = "Indianapolis Chapter of Indiana Black Expo, Inc."
$ws.Cells.Item(36,5).Value = 0
$ws.Cells.Item(36,6).Value = 0

$ws.Cells.Item(37,1).Value = 26
$ws.Cells.Item(37,2).Value = "OPERATIONAL"
$ws.Cells.Item(37,3).Value = "Indianapolis Competition Products"
$ws.Cells.Item(37,5).Value = 0
$ws.Cells.Item(37,6).Value = 0

$ws.Cells.Item(38,1).Value = 21
$ws.Cells.Item(38,2).Value = "OPERATIONAL"
$ws.Cells.Item(38,3).Value = "Indianapolis Motor Speedway"
$ws.Cells.Item(38,5).Value = 4.8
$ws.Cells.Item(38,6).Value = 11013

$ws.Cells.Item(39,1).Value = 13
$ws.Cells.Item(39,2).Value = "OPERATIONAL"
$ws.Cells.Item(39,3).Value = "Indy Air Expo"
$ws.Cells.Item(39,5).Value = 0
$ws.Cells.Item(39,6).Value = 0

$ws.Cells.Item(40,1).Value = 59
$ws.Cells.Item(40,2).Value = "OPERATIONAL"
$ws.Cells.Item(40,3).Value = "Ink Therapy Tattoo"
$ws.Cells.Item(40,5).Value = 4.7
$ws.Cells.Item(40,6).Value = 352

$ws.Cells.Item(41,1).Value = 58
$ws.Cells.Item(41,2).Value = "OPERATIONAL"
$ws.Cells.Item(41,3).Value = "Its A Block Party"
$ws.Cells.Item(41,5).Value = 4.7
$ws.Cells.Item(41,6).Value = 86

$ws.Cells.Item(42,1).Value = 10
$ws.Cells.Item(42,2).Value = "OPERATIONAL"
$ws.Cells.Item(42,3).Value = "Laser Storm"
$ws.Cells.Item(42,5).Value = 4
$ws.Cells.Item(42,6).Value = 51

$ws.Cells.Item(43,1).Value = 7
$ws.Cells.Item(43,2).Value = "OPERATIONAL"
$ws.Cells.Item(43,3).Value = "Nail Expo"
$ws.Cells.Item(43,5).Value = 3.2
$ws.Cells.Item(43,6).Value = 62

$ws.Cells.Item(44,1).Value = 8
$ws.Cells.Item(44,2).Value = "OPERATIONAL"
$ws.Cells.Item(44,3).Value = "National Expo, Inc"
$ws.Cells.Item(44,5).Value = 5
$ws.Cells.Item(44,6).Value = 1

$ws.Cells.Item(45,1).Value = 12
$ws.Cells.Item(45,2).Value = "OPERATIONAL"
$ws.Cells.Item(45,3).Value = "Off Road Expo"
$ws.Cells.Item(45,5).Value = 0
$ws.Cells.Item(45,6).Value = 0

$ws.Cells.Item(46,1).Value = 54
$ws.Cells.Item(46,2).Value = "OPERATIONAL"
$ws.Cells.Item(46,3).Value = "Pan Am Tower"
$ws.Cells.Item(46,5).Value = 4.3
$ws.Cells.Item(46,6).Value = 28

$ws.Cells.Item(47,1).Value = 35
$ws.Cells.Item(47,2).Value = "OPERATIONAL"
$ws.Cells.Item(47,3).Value = "Premier Surface"
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(47,6).Value = 0

$ws.Cells.Item(48,1).Value = 48
$ws.Cells.Item(48,2).Value = "OPERATIONAL"
$ws.Cells.Item(48,3).Value = "Purdue Extension / Horticulture Building"
$ws.Cells.Item(48,5).Value = 4.4
$ws.Cells.Item(48,6).Value = 36

$ws.Cells.Item(49,1).Value = 56
$ws.Cells.Item(49,2).Value = "OPERATIONAL"
$ws.Cells.Item(49,3).Value = "Roberts Camera"
$ws.Cells.Item(49,5).Value = 4.7
$ws.Cells.Item(49,6).Value = 518

$ws.Cells.Item(50,1).Value = 1
$ws.Cells.Item(50,2).Value = "OPERATIONAL"
$ws.Cells.Item(50,3).Value = "Royal Pin Expo"
$ws.Cells.Item(50,5).Value = 4.2
$ws.Cells.Item(50,6).Value = 1025

$ws.Cells.Item(51,1).Value = 31
$ws.Cells.Item(51,2).Value = "OPERATIONAL"
$ws.Cells.Item(51,3).Value = "Royal Pin Western"
$ws.Cells.Item(51,5).Value = 4.4
$ws.Cells.Item(51,6).Value = 1182

$ws.Cells.Item(52,1).Value = 16
$ws.Cells.Item(52,2).Value = "OPERATIONAL"
$ws.Cells.Item(52,3).Value = "Samps Hack Shack Brownsburg"
$ws.Cells.Item(52,5).Value = 4.9
$ws.Cells.Item(52,6).Value = 22

$ws.Cells.Item(53,1).Value = 32
$ws.Cells.Item(53,2).Value = "OPERATIONAL"
$ws.Cells.Item(53,3).Value = "Samps Hack Shack Plainfield"
$ws.Cells.Item(53,5).Value = 5
$ws.Cells.Item(53,6).Value = 9

$ws.Cells.Item(54,1).Value = 18
$ws.Cells.Item(54,2).Value = "OPERATIONAL"
$ws.Cells.Item(54,3).Value = "Shepard Events"
$ws.Cells.Item(54,5).Value = 3
$ws.Cells.Item(54,6).Value = 1

$ws.Cells.Item(55,1).Value = 36
$ws.Cells.Item(55,2).Value = "OPERATIONAL"
$ws.Cells.Item(55,3).Value = "Shepard Exposition Services"
$ws.Cells.Item(55,5).Value = 0
$ws.Cells.Item(55,6).Value = 0

$ws.Cells.Item(56,1).Value = 46
$ws.Cells.Item(56,2).Value = "OPERATIONAL"
$ws.Cells.Item(56,3).Value = "The Indiana Convention center"
$ws.Cells.Item(56,5).Value = 4.7
$ws.Cells.Item(56,6).Value = 21

$ws.Cells.Item(57,1).Value = 50
$ws.Cells.Item(57,2).Value = "OPERATIONAL"
$ws.Cells.Item(57,3).Value = "The Korner Garage"
$ws.Cells.Item(57,5).Value = 4.5
$ws.Cells.Item(57,6).Value = 6

$ws.Cells.Item(58,1).Value = 37
$ws.Cells.Item(58,2).Value = "OPERATIONAL"
$ws.Cells.Item(58,3).Value = "The Nest Event Center"
$ws.Cells.Item(58,5).Value = 4.8
$ws.Cells.Item(58,6).Value = 98

$ws.Cells.Item(59,1).Value = 39
$ws.Cells.Item(59,2).Value = "OPERATIONAL"
$ws.Cells.Item(59,3).Value = "West Pavilion"
$ws.Cells.Item(59,5).Value = 4.3
$ws.Cells.Item(59,6).Value = 33
